# "Refined files and redone analysis" -- update the computed enrichment /
# p-value table on Sheet1 with re-run numbers, and normalize the header /
# row-label cell formatting onto the already-present "text" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated analysis values (re-run statistics) ---------------------------
$ws.Range("B2").Value  = 2.4539514023930904
$ws.Range("C2").Value  = 0.026352068715589958
$ws.Range("D2").Value  = 3.5611369990680335
$ws.Range("E2").Value  = 0.00018907228507810821

$ws.Range("B3").Value  = 5.0204449450858917
$ws.Range("C3").Value  = 0.004433656180474612

$ws.Range("B6").Value  = 0.94221235664076952
$ws.Range("C6").Value  = 1.1384214525207235
$ws.Range("E6").Value  = 1.3513504984618245

$ws.Range("B8").Value  = 1.0921097770154373
$ws.Range("C8").Value  = 1.2076641939055082

$ws.Range("B9").Value  = 3.6039622641509435
$ws.Range("C9").Value  = 0.023137339199992207

$ws.Range("B10").Value = 0.7700774068698597
$ws.Range("C10").Value = 1.2108023052969186
$ws.Range("D10").Value = 1.2266773675762439
$ws.Range("E10").Value = 1.1464338447777489

# --- Re-normalize header row + row-label column styling ---------------------
# (re-applying the text format these cells already used collapses their
# style index back onto the workbook's canonical "text" cell style)
$ws.Range("A1:E1").NumberFormat = "@"
$ws.Range("A2:A10").NumberFormat = "@"

# --- Column width tweaks -----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.307291666666666
$ws.Columns.Item(5).ColumnWidth = 14.877604166666666
